$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.942.85'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '3.416.91'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '578.49'
$ws.Range("D6").Value = '144.51'
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.90%  '
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").Value = '4.002.72'
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '3.406.62'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '61.953.07'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").Value = '13.99'
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("E20").Value = '  +2.88%  '
$ws.Range("D21").Value = '391.28'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("D22").Value = '74.86'
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").Value = '0.553'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000115'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.557.11'
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("E27").Value = '  -2.10%  '
$ws.Range("E28").Value = '  +3.06%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '8.01'
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.19%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '23.56'
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("E35").Value = '  +6.02%  '
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '167.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("E38").Value = '  +4.45%  '
$ws.Range("D39").Value = '3.448.06'
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").Value = '28.83'
$ws.Range("E40").Value = '  +9.61%  '
$ws.Range("D41").Value = '0.0755'
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").Value = '0.785'
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("E44").Value = '  +1.38%  '
$ws.Range("D45").Value = '1.16'
$ws.Range("E45").Value = '  +4.15%  '
$ws.Range("D46").Value = '2.505.41'
$ws.Range("E46").Value = '  +2.10%  '
$ws.Range("D47").Value = '22.85'
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  -1.33%  '
